$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: comment on scanning data from input file -> new rubric text
$ws.Range("F29").Value = "(-1) for writing while inside do loop"
# Score for row 29 (scanning of data) changes from 8 to 15
$ws.Range("E29").Value = 15

# Row 30: comment on output -> new rubric text
$ws.Range("F30").Value = "(-4) for no output displayed due to compilation errors"

# Row 37: comment on compilation errors -> new rubric text
$ws.Range("F37").Value = "(-5) for compilation errors in CustomerMapping class"

# Move selection to F37 (also resets the scrolled top-left cell)
$ws.Range("F37").Select()
